$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Year of Treatment") entirely; this shifts C:K left to B:J
$ws.Range("B:B").Delete()

# Append ".jamais.jamais" to each header label in row 1 (B1:J1)
for ($col = 2; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $label = $cell.Value2
    $cell.Value = $label + ".jamais.jamais"
}
